$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append new row 37 ---
$ws = $wb.Worksheets.Item("Logs")

$row = 37
$ws.Cells.Item($row, 1).Value = "Terugbetaling"
$ws.Cells.Item($row, 2).Value = "support@testbedrijf123.nl"
$ws.Cells.Item($row, 3).Value = "Ik heb mijn retour gestuurd maar nog geen geld terug."
$ws.Cells.Item($row, 4).Value = "Intern verzoek / Actie voor medewerker"
$ws.Cells.Item($row, 5).Value = "Bedankt, we hebben dit doorgestuurd naar retour@testbedrijf123.nl."
$ws.Cells.Item($row, 6).Value = "2025-08-14 22:04:19"
$ws.Cells.Item($row, 7).Value = "Nee"
$ws.Cells.Item($row, 8).Value = "Ja"
$ws.Cells.Item($row, 9).Value = "Nee"
$ws.Cells.Item($row, 10).Value = "Nee"

# --- Extend conditional formatting ranges from row 36 to row 37 ---
foreach ($col in @("D","G","H","I","J")) {
    $oldRangeAddr = "$col" + "2:" + "$col" + "36"
    $newRangeAddr = "$col" + "2:" + "$col" + "37"
    $fc = $ws.Range($oldRangeAddr).FormatConditions
    for ($i = 1; $i -le $fc.Count; $i++) {
        $fc.Item($i).ModifyAppliesToRange($ws.Range($newRangeAddr))
    }
}

# --- Sheet "Dashboard": update count for "Intern verzoek / Actie voor medewerker" ---
$ws2 = $wb.Worksheets.Item("Dashboard")
$ws2.Cells.Item(2, 2).Value = 29
